$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows that are no longer part of the table (rows 6-12),
# shrinking the sheet down to the new A1:C5 extent.
$ws.Rows("6:12").Delete()

# Header row (Ballgorithm vs ESPN column headers) with new date label.
$ws.Range("A1").Value = "NBA, Monday 26th Feb 2024"
$ws.Range("B1").Value = "Ballgorithm"
$ws.Range("C1").Value = "ESPN"

# Fill column-by-column (matchups, then Ballgorithm picks, then ESPN picks)
# so new shared-string entries land in the same order the source workbook used.
$ws.Range("A2").Value = "Toronto Raptors (21-36) vs Indiana Pacers (33-26)"
$ws.Range("A3").Value = "Detroit Pistons (8-48) vs New York Knicks (34-23)"
$ws.Range("A4").Value = "Brooklyn Nets (21-35) vs Memphis Grizzlies (20-37)"
$ws.Range("A5").Value = "Miami Heat (31-25) vs Sacramento Kings (33-23)"

$ws.Range("B2").Value = "Indiana Pacers (63.33%)"
$ws.Range("B3").Value = "New York Knicks (67.86%)"
$ws.Range("B4").Value = "Memphis Grizzlies (56.67%)"
$ws.Range("B5").Value = "Sacramento Kings (64.00%)"

$ws.Range("C2").Value = "Indiana Pacers (65.7%)"
$ws.Range("C3").Value = "New York Knicks (76.8%)"
$ws.Range("C4").Value = "Brooklyn Nets (50.2%)"
$ws.Range("C5").Value = "Sacramento Kings (62.4%)"

# Match the saved selection from the authored workbook (active cell A5).
$ws.Range("A5").Select()
